$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the latest cryptos snapshot.
# D-column values are forced to text so Excel does not reinterpret
# dot-separated price strings (e.g. "26.082.45") as numbers/dates,
# then formatting is cleared so the cell keeps the workbook default style.

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.082.45"
$c.ClearFormats()
$ws.Range("E2").Value = "  -0.28%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.659.24"
$c.ClearFormats()
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -0.21%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "207.95"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.03%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5176"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  -0.22%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2581"
$c.ClearFormats()
$ws.Range("E8").Value = "  -3.40%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06285"
$c.ClearFormats()
$ws.Range("E9").Value = "  +0.14%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.90"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.77%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07530"
$c.ClearFormats()
$ws.Range("E11").Value = "  -0.04%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.666.04"
$c.ClearFormats()
$ws.Range("E12").Value = "  -0.65%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.401"
$c.ClearFormats()
$ws.Range("E13").Value = "  -1.38%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.5396"
$c.ClearFormats()
$ws.Range("E14").Value = "  -4.47%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "66.06"
$c.ClearFormats()
$ws.Range("E15").Value = "  -0.31%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.0₅7901"
$c.ClearFormats()
$ws.Range("E16").Value = "  -2.47%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "26.078.08"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("E18").Value = "  -0.14%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.690"
$c.ClearFormats()
$ws.Range("E19").Value = "  -3.04%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "187.74"
$c.ClearFormats()
$ws.Range("E20").Value = "  -0.09%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.17"
$c.ClearFormats()
$ws.Range("E21").Value = "  -3.04%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.185"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  -0.19%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "148.03"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.92%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.1206"
$c.ClearFormats()
$ws.Range("E25").Value = "  -3.89%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "7.370"
$c.ClearFormats()
$ws.Range("E26").Value = "  -3.32%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "15.61"
$c.ClearFormats()
$ws.Range("E27").Value = "  -1.67%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.383"
$c.ClearFormats()
$ws.Range("E28").Value = "  +2.90%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.06056"
$c.ClearFormats()
$ws.Range("E29").Value = "  -5.48%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.262"
$c.ClearFormats()
$ws.Range("E30").Value = "  -1.42%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.467"
$c.ClearFormats()
$ws.Range("E31").Value = "  -1.91%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.392"
$c.ClearFormats()
$ws.Range("E32").Value = "  -2.32%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.630"
$c.ClearFormats()
$ws.Range("E33").Value = "  -1.24%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.9825"
$c.ClearFormats()
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  -1.05%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.745"
$c.ClearFormats()
$ws.Range("E36").Value = "  +1.22%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.5860"
$c.ClearFormats()
$ws.Range("E37").Value = "  -3.51%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.102.93"
$c.ClearFormats()
$ws.Range("E38").Value = "  +0.29%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01592"
$c.ClearFormats()
$ws.Range("E39").Value = "  -1.14%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.947"
$c.ClearFormats()
$ws.Range("E40").Value = "  -3.38%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.8484"
$c.ClearFormats()
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("E42").Value = "  -0.34%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "99.84"
$c.ClearFormats()
$ws.Range("E43").Value = "  -0.14%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.811.56"
$c.ClearFormats()
$ws.Range("E44").Value = "  -0.96%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0₈110"
$c.ClearFormats()
$ws.Range("E45").Value = "  +1.20%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "55.02"
$c.ClearFormats()
$ws.Range("E46").Value = "  -3.01%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.9982"
$c.ClearFormats()
$ws.Range("E47").Value = "  -0.08%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.992"
$c.ClearFormats()
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  -0.66%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.4238"
$c.ClearFormats()
$ws.Range("E50").Value = "  -0.71%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "5.853"
$c.ClearFormats()
$ws.Range("E51").Value = "  -1.62%  "
